$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision drift on the existing A4 timestamp
$ws.Range("A4").Value = 45863.45853578704

# Append the new sensor reading row (row 5)
$ws.Range("A5").Value = 45863.54190494432
$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 19.44
$ws.Range("E5").Value = 74.63
$ws.Range("F5").Value = 611.72
$ws.Range("G5").Value = 10.73
$ws.Range("H5").Value = "ESE"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "13:00:20"

# Match the number format used by the other date/time cells in column A
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
